# Apply updated cryptocurrency price/volume data to Sheet1 (columns D and E, rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'23.316.08"
$ws.Range("E2").Value = "  -1.62%  "

$ws.Range("D3").Value = "'1.626.67"
$ws.Range("E3").Value = "  -1.65%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'1.003"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").Value = "'298.38"
$ws.Range("E6").Value = "  -1.55%  "

$ws.Range("D7").Value = "'0.3765"
$ws.Range("E7").Value = "  -0.99%  "

$ws.Range("D8").Value = "'50.13"
$ws.Range("E8").Value = "  -2.25%  "

$ws.Range("D9").Value = "'0.3515"
$ws.Range("E9").Value = "  -2.90%  "

$ws.Range("D10").Value = "'0.08011"
$ws.Range("E10").Value = "  -2.48%  "

$ws.Range("D11").Value = "'1.201"
$ws.Range("E11").Value = "  -3.55%  "

$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").Value = "'21.79"
$ws.Range("E13").Value = "  -3.74%  "

$ws.Range("D14").Value = "'6.257"
$ws.Range("E14").Value = "  -4.08%  "

$ws.Range("D15").Value = "'7.198"
$ws.Range("E15").Value = "  -3.24%  "

$ws.Range("D16").Value = "'0.00001186"
$ws.Range("E16").Value = "  -3.89%  "

$ws.Range("D17").Value = "'1.629.25"
$ws.Range("E17").Value = "  -1.54%  "

$ws.Range("D18").Value = "'95.17"
$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("D19").Value = "'0.06937"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("D20").Value = "'6.654"
$ws.Range("E20").Value = "  -2.24%  "

$ws.Range("D21").Value = "'17.21"
$ws.Range("E21").Value = "  -2.75%  "

$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").Value = "'12.21"
$ws.Range("E23").Value = "  -5.13%  "

$ws.Range("D24").Value = "'23.340.58"
$ws.Range("E24").Value = "  -1.56%  "

$ws.Range("D25").Value = "'2.446"
$ws.Range("E25").Value = "  -2.92%  "

$ws.Range("D26").Value = "'2.887"
$ws.Range("E26").Value = "  -5.22%  "

$ws.Range("D27").Value = "'20.74"
$ws.Range("E27").Value = "  -2.54%  "

$ws.Range("D28").Value = "'151.85"
$ws.Range("E28").Value = "  -1.04%  "

$ws.Range("D29").Value = "'5.168"
$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("D30").Value = "'131.78"
$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("D31").Value = "'1.813.20"
$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D32").Value = "'6.754"
$ws.Range("E32").Value = "  -2.57%  "

$ws.Range("D33").Value = "'2.117"
$ws.Range("E33").Value = "  -4.78%  "

$ws.Range("D34").Value = "'11.15"
$ws.Range("E34").Value = "  -4.71%  "

$ws.Range("D35").Value = "'0.9630"
$ws.Range("E35").Value = "  -9.86%  "

$ws.Range("D36").Value = "'0.02674"
$ws.Range("E36").Value = "  -4.99%  "

$ws.Range("D37").Value = "'0.08684"
$ws.Range("E37").Value = "  -1.31%  "

$ws.Range("D38").Value = "'0.2415"
$ws.Range("E38").Value = "  -4.66%  "

$ws.Range("D39").Value = "'5.796"
$ws.Range("E39").Value = "  -5.00%  "

$ws.Range("D40").Value = "'0.06745"
$ws.Range("E40").Value = "  -5.34%  "

$ws.Range("D41").Value = "'12.69"
$ws.Range("E41").Value = "  -2.12%  "

$ws.Range("D42").Value = "'0.6780"
$ws.Range("E42").Value = "  -3.54%  "

$ws.Range("D43").Value = "'1.293"
$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("E44").Value = "  -3.80%  "

$ws.Range("D45").Value = "'1.002"
$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("D46").Value = "'0.6268"
$ws.Range("E46").Value = "  -3.66%  "

$ws.Range("D47").Value = "'3.889"
$ws.Range("E47").Value = "  -2.27%  "

$ws.Range("D48").Value = "'2.220"
$ws.Range("E48").Value = "  -4.09%  "

$ws.Range("D49").Value = "'0.07653"
$ws.Range("E49").Value = "  -3.88%  "

$ws.Range("D50").Value = "'126.19"
$ws.Range("E50").Value = "  -1.48%  "

$ws.Range("D51").Value = "'1.196"
$ws.Range("E51").Value = "  -0.03%  "
